$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before row 8 (shifts old row 8 "Upload" and everything below down by one)
$ws.Rows.Item(8).Insert()

# Copy formatting from row 9 (the old row 8, "Upload"), range A9:P9, onto the new A8:P8
$ws.Range("A9:P9").Copy()
$ws.Range("A8:P8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new row 8 with "Force" and FALSE values
$ws.Cells.Item(8, 1).Value = "Force"
for ($col = 2; $col -le 16; $col++) {
    $ws.Cells.Item(8, $col).Value = $false
}
